$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Usuarios")

# Add three new user rows (14-16), matching the pattern of existing rows:
# column A = username, column B = email, column C = password ("Mynor123!")
$ws.Range("A14").Value = "Dani"
$ws.Range("B14").Value = "ellydanielabc@gmail.com"
$ws.Range("C14").Value = "Mynor123!"

$ws.Range("A15").Value = "Dani_"
$ws.Range("B15").Value = "dani@gmail.com"
$ws.Range("C15").Value = "Mynor123!"

$ws.Range("A16").Value = "Dani_1"
$ws.Range("B16").Value = "danibc@gmail.com"
$ws.Range("C16").Value = "Mynor123!"
